$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.218.93"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.477.75"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").Value = "2.477.56"
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D17").Value = "63.100.56"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").Value = "2.481.81"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("E22").Value = "  +9.28%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +16.76%  "
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "656.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.75%  "
$ws.Range("D28").Value = "0.0₃0990"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").Value = "2.607.47"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.84%  "
$ws.Range("E31").Value = "  +4.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("E35").Value = "  +4.23%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "0.0₆0323"
$ws.Range("E44").Value = "  -49.84%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "154.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0515"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
